$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).NumberFormat = "@"

$ws.Range("D2").Value = "25.973.73"
$ws.Range("E2").Value = "  +2.90%  "
$ws.Range("D3").Value = "1.600.80"
$ws.Range("E3").Value = "  +2.93%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "212.24"
$ws.Range("E5").Value = "  +2.69%  "
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("E7").Value = "  +1.11%  "
$ws.Range("E8").Value = "  +1.71%  "
$ws.Range("E9").Value = "  +0.60%  "
$ws.Range("D10").Value = "18.05"
$ws.Range("E10").Value = "  +1.52%  "
$ws.Range("E11").Value = "  +4.12%  "
$ws.Range("D12").Value = "1.823.50"
$ws.Range("E12").Value = "  +2.98%  "
$ws.Range("D13").Value = "1.593.31"
$ws.Range("E13").Value = "  +2.45%  "
$ws.Range("E14").Value = "  +0.14%  "
$ws.Range("E15").Value = "  +0.67%  "
$ws.Range("D16").Value = "25.973.61"
$ws.Range("E16").Value = "  +2.89%  "
$ws.Range("D17").Value = "60.24"
$ws.Range("E17").Value = "  +2.21%  "
$ws.Range("E18").Value = "  +1.92%  "
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("D20").Value = "201.31"
$ws.Range("E20").Value = "  +8.57%  "
$ws.Range("E21").Value = "  +2.54%  "
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("E23").Value = "  +2.59%  "
$ws.Range("D24").Value = "1.82"
$ws.Range("E24").Value = "  +10.40%  "
$ws.Range("D25").Value = "141.29"
$ws.Range("E25").Value = "  +0.33%  "
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("D27").Value = "0.122"
$ws.Range("E27").Value = "  -5.66%  "
$ws.Range("E28").Value = "  +1.87%  "
$ws.Range("E29").Value = "  +0.21%  "
$ws.Range("E30").Value = "  +1.74%  "
$ws.Range("D31").Value = "0.0472"
$ws.Range("E31").Value = "  +1.19%  "
$ws.Range("E32").Value = "  +2.20%  "
$ws.Range("E33").Value = "  -0.40%  "
$ws.Range("E34").Value = "  +1.43%  "
$ws.Range("E35").Value = "  +1.16%  "
$ws.Range("D36").Value = "0.0166"
$ws.Range("E36").Value = "  +11.59%  "
$ws.Range("D37").Value = "1.127.14"
$ws.Range("E37").Value = "  +3.96%  "
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("D39").Value = "0.791"
$ws.Range("E39").Value = "  +3.56%  "
$ws.Range("E40").Value = "  +2.53%  "
$ws.Range("D41").Value = "0.489"
$ws.Range("E41").Value = "  -0.86%  "
$ws.Range("E42").Value = "  -1.97%  "
$ws.Range("D43").Value = "1.735.62"
$ws.Range("E43").Value = "  +2.98%  "
$ws.Range("D44").Value = "5.12"
$ws.Range("E44").Value = "  +1.64%  "
$ws.Range("D45").Value = "93.04"
$ws.Range("E45").Value = "  +0.31%  "
$ws.Range("E46").Value = "  +3.62%  "
$ws.Range("D47").Value = "53.27"
$ws.Range("E47").Value = "  +1.87%  "
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("D49").Value = "0.408"
$ws.Range("E49").Value = "  +1.03%  "
$ws.Range("D50").Value = "1.01"
$ws.Range("E50").Value = "  +0.37%  "
$ws.Range("D51").Value = "0.0₇0924"
$ws.Range("E51").Value = "  -16.96%  "
